$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell reference -> new text value, derived from the commit diff.
$updates = @{
    'D2' = '64.013.82'
    'E2' = '  -0.43%  '
    'D3' = '3.468.24'
    'E3' = '  -0.70%  '
    'E4' = '  +0.07%  '
    'D5' = '583.61'
    'E5' = '  -0.63%  '
    'E6' = '  -2.67%  '
    'E7' = '  +0.02%  '
    'D8' = '0.481'
    'E8' = '  -1.17%  '
    'D9' = '7.66'
    'E9' = '  +5.52%  '
    'E10' = '  -1.09%  '
    'E11' = '  -0.41%  '
    'D12' = '4.064.38'
    'E12' = '  -0.52%  '
    'E13' = '  -0.16%  '
    'E14' = '  -2.60%  '
    'D15' = '3.472.00'
    'E15' = '  -0.99%  '
    'D16' = '63.999.62'
    'E16' = '  -0.56%  '
    'D17' = '24.90'
    'E17' = '  -3.41%  '
    'D18' = '9.93'
    'E18' = '  +0.44%  '
    'D19' = '5.66'
    'E19' = '  -1.63%  '
    'D20' = '13.35'
    'E20' = '  -1.89%  '
    'D21' = '383.88'
    'E21' = '  -2.67%  '
    'D22' = '0.564'
    'E22' = '  -1.11%  '
    'D23' = '3.609.18'
    'E23' = '  -0.63%  '
    'D24' = '74.39'
    'E24' = '  -0.48%  '
    'D25' = '1.00'
    'E25' = '  +0.09%  '
    'D26' = '5.38'
    'E26' = '  -6.24%  '
    'E27' = '  -3.68%  '
    'E28' = '  +1.50%  '
    'D29' = '2.21'
    'E29' = '  -0.58%  '
    'D30' = '7.04'
    'E30' = '  -4.21%  '
    'E31' = '  +3.57%  '
    'D32' = '7.92'
    'E32' = '  -3.64%  '
    'D33' = '1.42'
    'D34' = '3.497.67'
    'E34' = '  -0.49%  '
    'D36' = '22.87'
    'E36' = '  -2.28%  '
    'E37' = '  +0.92%  '
    'D38' = '6.74'
    'E38' = '  -2.20%  '
    'D39' = '163.32'
    'E39' = '  -1.86%  '
    'E40' = '  -3.92%  '
    'E41' = '  -1.02%  '
    'E42' = '  -1.32%  '
    'E43' = '  +0.07%  '
    'D44' = '41.49'
    'E44' = '  -0.81%  '
    'D45' = '4.33'
    'E45' = '  -1.17%  '
    'E46' = '  -2.32%  '
    'D47' = '23.52'
    'E47' = '  -6.72%  '
    'D48' = '1.11'
    'E48' = '  -4.03%  '
    'D49' = '6.69'
    'E49' = '  -1.00%  '
    'D50' = '0.891'
    'E50' = '  -0.11%  '
    'D51' = '2.320.54'
    'E51' = '  -5.61%  '
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $cell = $ws.Range($cellRef)

    # Force the cell to be treated as plain text so values such as
    # "24.90", "1.00" or "64.013.82" are not re-interpreted/rounded as
    # numbers by Excel. We restore the default "Normal" style afterwards
    # so no stray formatting is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}
